$d = $word.ActiveDocument

# 1. Add a new "INVOICE NO: 8668" paragraph right after the GSTIN paragraph,
#    inheriting the same (bold, sz 24) paragraph/run formatting.
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "GSTIN: 33BESPV7542D1Z3`r") {
        $p.Range.InsertParagraphAfter()
        $newP = $d.Paragraphs.Item($i + 1)
        $newP.Range.Text = "INVOICE NO: 8668"
        break
    }
}

# 2. Merge the split "Total # of Cabs: " + "71" runs into a single run.
$d.Content.Find.Execute("Total # of Cabs: 71", $true, $false, $false, $false, $false, $true, 1, $false, "Total # of Cabs: 71", 2) | Out-Null

# 3. Merge the split "Amount: Rs. 5" + "," + "63" + "," + "992/-" runs into a single run.
$d.Content.Find.Execute("Amount: Rs. 5,63,992/-", $true, $false, $false, $false, $false, $true, 1, $false, "Amount: Rs. 5,63,992/-", 2) | Out-Null
